$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.702.28'
$ws.Range('E2').Value = '  -8.79%  '
$ws.Range('D3').Value = '2.880.64'
$ws.Range('E3').Value = '  -8.86%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '547.47'
$ws.Range('E5').Value = '  -9.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '120.45'
$ws.Range('E6').Value = '  -11.01%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').Value = '2.876.51'
$ws.Range('E8').Value = '  -8.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('E9').Value = '  -4.21%  '
$ws.Range('E10').Value = '  -12.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.79'
$ws.Range('E11').Value = '  -10.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.428'
$ws.Range('E12').Value = '  -5.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000211'
$ws.Range('E13').Value = '  -11.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '31.16'
$ws.Range('E14').Value = '  -10.44%  '
$ws.Range('E15').Value = '  -1.58%  '
$ws.Range('D16').Value = '3.356.40'
$ws.Range('E16').Value = '  -8.79%  '
$ws.Range('D17').Value = '2.884.52'
$ws.Range('E17').Value = '  -8.88%  '
$ws.Range('D18').Value = '57.932.90'
$ws.Range('E18').Value = '  -8.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.30'
$ws.Range('E19').Value = '  -3.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '412.61'
$ws.Range('E20').Value = '  -10.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.71'
$ws.Range('E21').Value = '  -8.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.646'
$ws.Range('E22').Value = '  -7.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.74'
$ws.Range('E23').Value = '  -11.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.46'
$ws.Range('E24').Value = '  -6.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '76.31'
$ws.Range('E25').Value = '  -8.18%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.43'
$ws.Range('E28').Value = '  -9.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.89'
$ws.Range('E29').Value = '  -8.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.98'
$ws.Range('E30').Value = '  -9.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '24.46'
$ws.Range('E31').Value = '  -9.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.89'
$ws.Range('E32').Value = '  -12.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0925'
$ws.Range('E33').Value = '  -7.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.32'
$ws.Range('E34').Value = '  -9.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '48.43'
$ws.Range('E35').Value = '  -5.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.882'
$ws.Range('E36').Value = '  -13.18%  '
$ws.Range('E37').Value = '  -17.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.31'
$ws.Range('E38').Value = '  +2.41%  '
$ws.Range('D39').Value = '0.0₃0615'
$ws.Range('E39').Value = '  -15.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0340'
$ws.Range('E40').Value = '  -12.77%  '
$ws.Range('E41').Value = '  -7.83%  '
$ws.Range('D42').Value = '2.600.71'
$ws.Range('E42').Value = '  -6.69%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.31'
$ws.Range('E44').Value = '  -11.35%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '345.45'
$ws.Range('E45').Value = '  -11.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '117.86'
$ws.Range('E46').Value = '  -6.15%  '
$ws.Range('E47').Value = '  -9.53%  '
$ws.Range('E48').Value = '  -5.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.91'
$ws.Range('E49').Value = '  -9.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.44'
$ws.Range('E50').Value = '  -10.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.92'
$ws.Range('E51').Value = '  -10.79%  '
